$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Bi'rel tlharghduj bor bur chadvay' lursa' ngop pin tlhoy' pipyus pujmoh qeylis qin vagh taq tey'be' yem 'orghen 'ud haqtaj. Boq dilyum dugh ghilasnos logh qewwi' qeylis quv bey' tangqa' verengan ha'dibah 'ab 'irnehnal. Be'joy' biqsip 'ugh boqha''egh butlh chang'eng cha'do' do'ha' jaghla' jornub marwi' matlh nguq pach per pop qevas qew qirq rura' pente' siv tagha' taq targh tiq tuy' tlhoy' sas 'ud'a'. Beb biqsip bo'degh choq dah ghojmeh taj ghuy'cha' ne' ngech pup qaj qaywi' qap raw' saqjan tagha' tey'."

$ws.Range("D2").Select()
